# Auto-generated edit script applying value changes per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 2.88
$ws.Range("H2").Value = 3.1
$ws.Range("I2").Value = 2.55
$ws.Range("P2").Value = 1.41
$ws.Range("Q2").Value = 2.62
$ws.Range("X2").Value = 23
$ws.Range("AH2").Value = 26

# Row 3
$ws.Range("N3").Value = 1.84
$ws.Range("O3").Value = 2.06
$ws.Range("P3").Value = 1.3

# Row 4
$ws.Range("J4").Value = 1.02
$ws.Range("L4").Value = 1.15
$ws.Range("M4").Value = 4.5
$ws.Range("N4").Value = 1.6
$ws.Range("O4").Value = 2.3

# Row 5
$ws.Range("J5").Value = 1.03
$ws.Range("L5").Value = 1.27
$ws.Range("N5").Value = 1.92
$ws.Range("O5").Value = 1.77

# Row 6
$ws.Range("J6").Value = 1.1
$ws.Range("L6").Value = 1.58

# Row 11
$ws.Range("K11").Value = 8

# Row 19
$ws.Range("G19").Value = 2.77
$ws.Range("H19").Value = 3.45
$ws.Range("I19").Value = 2.4
$ws.Range("K19").Value = 7.5
$ws.Range("L19").Value = 1.34
$ws.Range("M19").Value = 3.1
$ws.Range("N19").Value = 2
$ws.Range("P19").Value = 1.39
$ws.Range("Q19").Value = 2.85
$ws.Range("R19").Value = 1.8
$ws.Range("T19").Value = 8.25
$ws.Range("U19").Value = 14.5
$ws.Range("V19").Value = 11
$ws.Range("W19").Value = 35
$ws.Range("X19").Value = 26
$ws.Range("Z19").Value = 7.5
$ws.Range("AB19").Value = 16
$ws.Range("AD19").Value = 700
$ws.Range("AE19").Value = 7.6
$ws.Range("AF19").Value = 12
$ws.Range("AG19").Value = 10
$ws.Range("AH19").Value = 26
$ws.Range("AI19").Value = 22

# Row 22
$ws.Range("J22").Value = 1.07
$ws.Range("K22").Value = 9
$ws.Range("N22").Value = 2.1
$ws.Range("O22").Value = 1.7

# Row 24
$ws.Range("N24").Value = 2
$ws.Range("O24").Value = 1.8

# Row 25
$ws.Range("N25").Value = 1.57

# Row 27
$ws.Range("G27").Value = 4.85
$ws.Range("I27").Value = 1.65
$ws.Range("T27").Value = 8.75
$ws.Range("U27").Value = 20
$ws.Range("AE27").Value = 4.7
$ws.Range("AG27").Value = 7.2
$ws.Range("AH27").Value = 9.75

# Row 28
$ws.Range("G28").Value = 3.6
$ws.Range("H28").Value = 3.2
$ws.Range("I28").Value = 1.91
$ws.Range("N28").Value = 1.91
$ws.Range("O28").Value = 1.7
$ws.Range("T28").Value = 8.75
$ws.Range("U28").Value = 16
$ws.Range("V28").Value = 10.25
$ws.Range("X28").Value = 27
$ws.Range("Y28").Value = 30
$ws.Range("Z28").Value = 8.75
$ws.Range("AA28").Value = 5.5
$ws.Range("AB28").Value = 12
$ws.Range("AC28").Value = 50
$ws.Range("AD28").Value = 350
$ws.Range("AE28").Value = 5.9
$ws.Range("AF28").Value = 7.6
$ws.Range("AG28").Value = 7.1
$ws.Range("AH28").Value = 13.5
$ws.Range("AJ28").Value = 22

# Row 30
$ws.Range("K30").Value = 19

# Row 34
$ws.Range("L34").Value = 1.17
$ws.Range("M34").Value = 5
$ws.Range("N34").Value = 1.57
$ws.Range("O34").Value = 2.35

# Row 35
$ws.Range("G35").Value = 7
$ws.Range("H35").Value = 5
$ws.Range("I35").Value = 1.36
$ws.Range("K35").Value = 17
$ws.Range("N35").Value = 1.57
$ws.Range("O35").Value = 2.3
$ws.Range("U35").Value = 41
$ws.Range("V35").Value = 21
$ws.Range("W35").Value = 81
$ws.Range("X35").Value = 51
$ws.Range("Y35").Value = 51
$ws.Range("AA35").Value = 10
$ws.Range("AB35").Value = 21
$ws.Range("AD35").Value = 301
$ws.Range("AF35").Value = 7
$ws.Range("AH35").Value = 9

# Row 37
$ws.Range("N37").Value = 1.88
$ws.Range("O37").Value = 1.98

# Row 39
$ws.Range("N39").Value = 2.87
$ws.Range("O39").Value = 1.37

# Row 42
$ws.Range("G42").Value = 1.14
$ws.Range("H42").Value = 10
$ws.Range("I42").Value = 12
$ws.Range("K42").Value = 17
$ws.Range("N42").Value = 1.25
$ws.Range("O42").Value = 3.75
$ws.Range("R42").Value = 1.83
$ws.Range("S42").Value = 1.83
$ws.Range("U42").Value = 9
$ws.Range("AB42").Value = 29
$ws.Range("AC42").Value = 51
$ws.Range("AD42").Value = 500
$ws.Range("AG42").Value = 29
$ws.Range("AI42").Value = 67

# Row 44
$ws.Range("H44").Value = 3.3
$ws.Range("I44").Value = 2.1
$ws.Range("L44").Value = 1.25
$ws.Range("M44").Value = 3.2
$ws.Range("N44").Value = 1.75
$ws.Range("O44").Value = 1.85
$ws.Range("R44").Value = 1.62
$ws.Range("S44").Value = 2.02
$ws.Range("T44").Value = 11.25
$ws.Range("U44").Value = 19
$ws.Range("V44").Value = 11
$ws.Range("X44").Value = 26
$ws.Range("Y44").Value = 30
$ws.Range("Z44").Value = 10.75
$ws.Range("AA44").Value = 6.4
$ws.Range("AB44").Value = 12.5
$ws.Range("AC44").Value = 50
$ws.Range("AD44").Value = 350
$ws.Range("AE44").Value = 8
$ws.Range("AF44").Value = 10.5
$ws.Range("AG44").Value = 8.5
$ws.Range("AH44").Value = 20
$ws.Range("AI44").Value = 16.5
$ws.Range("AJ44").Value = 25

# Row 45
$ws.Range("G45").Value = 2.55
$ws.Range("H45").Value = 3.6
$ws.Range("I45").Value = 2.37
$ws.Range("L45").Value = 1.22
$ws.Range("M45").Value = 3.45
$ws.Range("N45").Value = 1.65
$ws.Range("O45").Value = 1.98
$ws.Range("R45").Value = 1.57
$ws.Range("S45").Value = 2.1
$ws.Range("T45").Value = 10.25
$ws.Range("U45").Value = 14
$ws.Range("X45").Value = 19.5
$ws.Range("Y45").Value = 26
$ws.Range("Z45").Value = 13
$ws.Range("AA45").Value = 7.1
$ws.Range("AB45").Value = 13
$ws.Range("AC45").Value = 50
$ws.Range("AD45").Value = 350
$ws.Range("AE45").Value = 9.75
$ws.Range("AF45").Value = 13
$ws.Range("AI45").Value = 18
$ws.Range("AJ45").Value = 25

# Row 48
$ws.Range("G48").Value = 1.78
$ws.Range("H48").Value = 3.55
$ws.Range("I48").Value = 4.05
$ws.Range("L48").Value = 1.27
$ws.Range("M48").Value = 3.1
$ws.Range("N48").Value = 1.78
$ws.Range("O48").Value = 1.82
$ws.Range("P48").Value = 1.4
$ws.Range("Q48").Value = 2.52
$ws.Range("R48").Value = 1.72
$ws.Range("S48").Value = 1.88
$ws.Range("T48").Value = 7.2
$ws.Range("U48").Value = 8.5
$ws.Range("V48").Value = 8.25
$ws.Range("X48").Value = 14
$ws.Range("Y48").Value = 25
$ws.Range("Z48").Value = 10.25
$ws.Range("AA48").Value = 6.9
$ws.Range("AB48").Value = 15
$ws.Range("AC48").Value = 70
$ws.Range("AD48").Value = 500
$ws.Range("AE48").Value = 11.75
$ws.Range("AF48").Value = 23
$ws.Range("AH48").Value = 65
$ws.Range("AI48").Value = 37
$ws.Range("AJ48").Value = 45

# Row 53
$ws.Range("I53").Value = 3.25
$ws.Range("L53").Value = 1.28
$ws.Range("T53").Value = 8
$ws.Range("U53").Value = 10.75
$ws.Range("Y53").Value = 25
$ws.Range("AE53").Value = 11
$ws.Range("AF53").Value = 18.5
$ws.Range("AI53").Value = 27
$ws.Range("AJ53").Value = 30

# Row 55
$ws.Range("J55").Value = 1.07
$ws.Range("K55").Value = 9
$ws.Range("N55").Value = 2.15
$ws.Range("O55").Value = 1.67

# Row 58
$ws.Range("L58").Value = 1.25
$ws.Range("M58").Value = 3.75
$ws.Range("N58").Value = 1.88
$ws.Range("O58").Value = 1.93

# Row 59
$ws.Range("G59").Value = 3.15
$ws.Range("H59").Value = 2.67
$ws.Range("I59").Value = 2.5
$ws.Range("J59").Value = 1.14
$ws.Range("K59").Value = 4.3
$ws.Range("L59").Value = 1.6
$ws.Range("M59").Value = 2.07
$ws.Range("N59").Value = 2.67
$ws.Range("O59").Value = 1.36
$ws.Range("P59").Value = 1.62
$ws.Range("Q59").Value = 2.05
$ws.Range("R59").Value = 2.15
$ws.Range("T59").Value = 6.4
$ws.Range("U59").Value = 14
$ws.Range("W59").Value = 45
$ws.Range("X59").Value = 40
$ws.Range("Y59").Value = 65
$ws.Range("Z59").Value = 4.6
$ws.Range("AA59").Value = 5.5
$ws.Range("AB59").Value = 19.5
$ws.Range("AE59").Value = 5.8
$ws.Range("AF59").Value = 10.75
$ws.Range("AG59").Value = 10.25
$ws.Range("AH59").Value = 29
$ws.Range("AI59").Value = 27
$ws.Range("AJ59").Value = 50
